# feat: add 2022-Q1 data
#
# The workbook originally has 3 sheets: 2021-Q3, 2021-Q4, 总计.
# We need to end up with 4 sheets: 2021-Q3, 2021-Q4, 2022-Q1, 总计
#   - "2022-Q1" (new) holds the per-fund holdings snapshot for the new quarter.
#   - "总计" keeps its original per-quarter summary rows, with a new row
#     inserted at the top for 2022-Q1 (4 funds, 0.33 亿元).
#
# To land the new sheet with the same sheetId/rId sequencing as the target
# (2022-Q1 = sheetId 3 / rId3, 总计 = sheetId 4 / rId4), we rename the
# existing "总计" sheet (sheetId 3) to "2022-Q1" and rewrite its data, then
# add a brand-new sheet right after it and name that one "总计".

$wb = $excel.ActiveWorkbook
$wb.Application.DisplayAlerts = $false

# ---------------------------------------------------------------------
# Step 1: turn the existing "总计" sheet into the new "2022-Q1" sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Extend the header/index-column styling (the style index carried by D1 /
# A3 -- bold + bordered + centered) onto the new cells we are about to
# fill in, so formatting stays consistent with the rest of the sheet.
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q1.Range("A3").Copy()
$q1.Range("A4:A5").PasteSpecial(-4122)   # xlPasteFormats
$wb.Application.CutCopyMode = $false

# Header row (plain text, not numeric-looking -> safe to set directly)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Formula = '="000800"'
$q1.Range("C2").Value = "华商未来主题混合"
$q1.Range("D2").Formula = '="4.31"'
$q1.Range("E2").Formula = '="84.71"'
$q1.Range("F2").Formula = '="3.71"'
$q1.Range("G2").Formula = '="0.1599"'
$q1.Range("H2").Value = 6

# Row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Formula = '="410001"'
$q1.Range("C3").Value = "华富竞争力优选混合"
$q1.Range("D3").Formula = '="3.34"'
$q1.Range("E3").Formula = '="89.23"'
$q1.Range("F3").Formula = '="3.84"'
$q1.Range("G3").Formula = '="0.1283"'
$q1.Range("H3").Value = 10

# Row 4 (new)
$q1.Range("A4").Value = 2
$q1.Range("B4").Formula = '="007713"'
$q1.Range("C4").Value = "华富科技动能混合"
$q1.Range("D4").Formula = '="0.56"'
$q1.Range("E4").Formula = '="86.98"'
$q1.Range("F4").Formula = '="4.31"'
$q1.Range("G4").Formula = '="0.0241"'
$q1.Range("H4").Value = 10

# Row 5 (new)
$q1.Range("A5").Value = 3
$q1.Range("B5").Formula = '="009999"'
$q1.Range("C5").Value = "东方中国红利混合"
$q1.Range("D5").Formula = '="0.69"'
$q1.Range("E5").Formula = '="72.52"'
$q1.Range("F5").Formula = '="2.79"'
$q1.Range("G5").Formula = '="0.0193"'
$q1.Range("H5").Value = 5

# The cells above that carry numeric-looking text (fund codes, percentages,
# …) were written as `="literal"` formulas so Excel wouldn't silently
# reinterpret the strings as numbers. Freeze them down to plain text values
# now (one bulk copy/paste-values over the whole block covers the plain
# literals too, which is a harmless no-op for them).
$q1.Range("B2:G5").Copy()
$q1.Range("B2:G5").PasteSpecial(-4163)   # xlPasteValues
$wb.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 2: add a fresh "总计" sheet right after "2022-Q1", reproducing the
# original summary table plus the new 2022-Q1 row on top.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Borrow formatting from the "2022-Q1" sheet (previously "总计") for the
# header row and the index column, so the new sheet matches the original
# look (bold header / bold bordered index column).
$q1.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)   # xlPasteFormats

$q1.Range("A2:A4").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)   # xlPasteFormats
$total.Application.CutCopyMode = $false

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.33

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.02

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 6
$total.Range("D4").Value = 0.75
